$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrección header: Description de la fila 2 (B2="61200821621") cambia de
# "Alimentacion, HeaderLEDs, LCD" a "Alimentacion, LCD"
# (se usa comilla inicial para conservar el quotePrefix / estilo original de la celda)
$ws.Range("D2").Value = "'Alimentacion, LCD"

# Corrección cantidad y subtotal (fila 2)
$ws.Range("H2").Value = 2
$ws.Range("O2").Value = 1.94

# Corrección capacitores (fila 8)
$ws.Range("N8").Value = 0.17979000000000001

# Corrección resistencias (fila 14)
$ws.Range("N14").Value = 0.089700000000000002
$ws.Range("O14").Value = 0.089700000000000002
